$d = $word.ActiveDocument

# "... jumlah 18 pemain, ..." -> "... jumlah 19 pemain, ..."
$rng1 = $d.Content.Duplicate
$rng1.Find.Execute("jumlah 18 pemain", $true, $false, $false, $false, $false,
                    $true, 1, $false, "jumlah 19 pemain", 2) | Out-Null

# "... didesuaikan kembali." -> "... disesuaikan kembali."
$rng2 = $d.Content.Duplicate
$rng2.Find.Execute("didesuaikan kembali", $true, $false, $false, $false, $false,
                    $true, 1, $false, "disesuaikan kembali", 2) | Out-Null
